$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Nombre" (count) column ------------------------------
# This shifts former columns C:F (Capacité, Taille, Consommation, Frais)
# left into B:E.
$ws.Columns("B").Delete()

# --- Insert extra rows so each vehicle model gets one row per unit ---
# After the column delete the sheet looks like:
#   row1 header
#   row2 Poids Lourd
#   row3 Petit Poids Lourd
#   row4 Camion Frigo   (only one row so far, need 7 total -> insert 6)
#   row5 Fourgon        (only one row so far, need 3 total -> insert 2)
$ws.Rows("5:10").Insert()   # 6 new blank rows after the Camion Frigo row (now rows 5-10)
$ws.Rows("11:12").Insert()  # 2 new blank rows after the Fourgon row (now rows 11-12)

# --- Rewrite the whole table with the final values --------------------
$headers = @("Nom", "Capacité (kg)", "Taille(Palettes)", "Consommation (L/100km)", "Frais")
for ($col = 1; $col -le 5; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$data = @(
    @("Poids Lourd", 19000, 12, 35, "Non"),
    @("Petit Poids Lourd", 6000, 6, 23, "Non"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Camion Frigo", 900, 4, 8, "Oui"),
    @("Fourgon", 1200, 3, 8, "Non"),
    @("Fourgon", 1200, 3, 8, "Non"),
    @("Fourgon", 1200, 3, 8, "Non")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# --- Cosmetic bits that mirror the authored workbook -------------------
$ws.Range("D7").Select()
